$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the cells can be edited
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text
$disclaimerCell = $ws.Range("A9")
$disclaimerCell.Value = $disclaimerCell.Value2 -replace "2021-03-23", "2021-03-24"

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2467879585299595
$ws.Range("E2").Value = 0.007227332457293123

$ws.Range("D3").Value = 0.2451895513056127
$ws.Range("E3").Value = 0.00361881785283491

$ws.Range("D4").Value = 0.2519872516776257
$ws.Range("E4").Value = -0.01213277375047705

$ws.Range("D5").Value = 0.256035238486802
$ws.Range("E5").Value = -0.0252156602521566

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = -0.00684248695052625

# Restore sheet protection to its original (protected) state
$ws.Protect()
